# Apply weekly update: for rows 2..19 the Volumen/Precio/Fecha block
# (columns D, M, N, O, P, S) gets re-shuffled between rows according to
# a fixed permutation (source row for each destination row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (i.e. new value at $dest = old value at $src)
$map = @{
    2  = 16
    3  = 6
    4  = 15
    5  = 10
    6  = 2
    7  = 9
    8  = 11
    9  = 12
    10 = 8
    11 = 18
    12 = 17
    13 = 7
    14 = 19
    15 = 4
    16 = 13
    17 = 14
    18 = 3
    19 = 5
}

$cols = @("D", "M", "N", "O", "P", "S")

# Snapshot the original values for the affected columns/rows before
# writing anything, since several cells are both a source and a
# destination in the permutation.
$orig = @{}
foreach ($row in 2..19) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $orig[$addr] = $ws.Range($addr).Value2
    }
}

foreach ($row in 2..19) {
    $src = $map[$row]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $orig["$col$src"]
    }
}
